$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "72+16="
$t.Cell(1, 2).Range.Text = "67+7="
$t.Cell(1, 3).Range.Text = "83-64="
$t.Cell(1, 4).Range.Text = "25+4="
$t.Cell(1, 5).Range.Text = "3+51="
$t.Cell(2, 1).Range.Text = "54+20="
$t.Cell(2, 2).Range.Text = "17+21="
$t.Cell(2, 3).Range.Text = "91-79="
$t.Cell(2, 4).Range.Text = "68-57="
$t.Cell(2, 5).Range.Text = "39-17="
$t.Cell(3, 1).Range.Text = "12+53="
$t.Cell(3, 2).Range.Text = "71-30="
$t.Cell(3, 3).Range.Text = "14+68="
$t.Cell(3, 4).Range.Text = "44+21="
$t.Cell(3, 5).Range.Text = "59-0="
$t.Cell(4, 1).Range.Text = "83-41="
$t.Cell(4, 2).Range.Text = "77+9="
$t.Cell(4, 3).Range.Text = "27+52="
$t.Cell(4, 4).Range.Text = "84-50="
$t.Cell(4, 5).Range.Text = "47-21="
$t.Cell(5, 1).Range.Text = "9+15="
$t.Cell(5, 2).Range.Text = "55-29="
$t.Cell(5, 3).Range.Text = "12+6="
$t.Cell(5, 4).Range.Text = "40+55="
$t.Cell(5, 5).Range.Text = "13+77="
$t.Cell(6, 1).Range.Text = "12+58="
$t.Cell(6, 2).Range.Text = "47+16="
$t.Cell(6, 3).Range.Text = "64-10="
$t.Cell(6, 4).Range.Text = "11+5="
$t.Cell(6, 5).Range.Text = "90-69="
$t.Cell(7, 1).Range.Text = "25-1="
$t.Cell(7, 2).Range.Text = "44+27="
$t.Cell(7, 3).Range.Text = "67-35="
$t.Cell(7, 4).Range.Text = "9+83="
$t.Cell(7, 5).Range.Text = "88+8="
$t.Cell(8, 1).Range.Text = "64-61="
$t.Cell(8, 2).Range.Text = "48-7="
$t.Cell(8, 3).Range.Text = "43-37="
$t.Cell(8, 4).Range.Text = "43+11="
$t.Cell(8, 5).Range.Text = "19+76="
$t.Cell(9, 1).Range.Text = "33+51="
$t.Cell(9, 2).Range.Text = "22+74="
$t.Cell(9, 3).Range.Text = "21+65="
$t.Cell(9, 4).Range.Text = "42-41="
$t.Cell(9, 5).Range.Text = "84-53="
$t.Cell(10, 1).Range.Text = "14+84="
$t.Cell(10, 2).Range.Text = "53-44="
$t.Cell(10, 3).Range.Text = "39+35="
$t.Cell(10, 4).Range.Text = "88-13="
$t.Cell(10, 5).Range.Text = "87-77="
$t.Cell(11, 1).Range.Text = "39+26="
$t.Cell(11, 2).Range.Text = "37-9="
$t.Cell(11, 3).Range.Text = "54+4="
$t.Cell(11, 4).Range.Text = "96-91="
$t.Cell(11, 5).Range.Text = "23-14="
$t.Cell(12, 1).Range.Text = "73-30="
$t.Cell(12, 2).Range.Text = "61+24="
$t.Cell(12, 3).Range.Text = "19+34="
$t.Cell(12, 4).Range.Text = "19-6="
$t.Cell(12, 5).Range.Text = "34+47="
$t.Cell(13, 1).Range.Text = "84-13="
$t.Cell(13, 2).Range.Text = "35+15="
$t.Cell(13, 3).Range.Text = "41-27="
$t.Cell(13, 4).Range.Text = "97-25="
$t.Cell(13, 5).Range.Text = "42+35="
$t.Cell(14, 1).Range.Text = "3+57="
$t.Cell(14, 2).Range.Text = "54-29="
$t.Cell(14, 3).Range.Text = "75-48="
$t.Cell(14, 4).Range.Text = "80-29="
$t.Cell(14, 5).Range.Text = "10+78="
$t.Cell(15, 1).Range.Text = "91-39="
$t.Cell(15, 2).Range.Text = "35+23="
$t.Cell(15, 3).Range.Text = "58-12="
$t.Cell(15, 4).Range.Text = "73-67="
$t.Cell(15, 5).Range.Text = "87-81="
$t.Cell(16, 1).Range.Text = "8+72="
$t.Cell(16, 2).Range.Text = "28+16="
$t.Cell(16, 3).Range.Text = "73+4="
$t.Cell(16, 4).Range.Text = "96-56="
$t.Cell(16, 5).Range.Text = "81-12="
$t.Cell(17, 1).Range.Text = "50+37="
$t.Cell(17, 2).Range.Text = "52-5="
$t.Cell(17, 3).Range.Text = "46-21="
$t.Cell(17, 4).Range.Text = "16-10="
$t.Cell(17, 5).Range.Text = "42+24="
$t.Cell(18, 1).Range.Text = "87-46="
$t.Cell(18, 2).Range.Text = "27+45="
$t.Cell(18, 3).Range.Text = "91-49="
$t.Cell(18, 4).Range.Text = "83-70="
$t.Cell(18, 5).Range.Text = "77-25="
$t.Cell(19, 1).Range.Text = "86+10="
$t.Cell(19, 2).Range.Text = "18+29="
$t.Cell(19, 3).Range.Text = "67-41="
$t.Cell(19, 4).Range.Text = "41-18="
$t.Cell(19, 5).Range.Text = "42-5="
$t.Cell(20, 1).Range.Text = "56+43="
$t.Cell(20, 2).Range.Text = "26+24="
$t.Cell(20, 3).Range.Text = "54+22="
$t.Cell(20, 4).Range.Text = "16+80="
$t.Cell(20, 5).Range.Text = "1+49="
